$d = $word.ActiveDocument

$d.Content.Find.Execute("136÷6=22, 4", $true, $false, $false, $false, $false, $true, 1, $false, "307÷3=102, 1", 2) | Out-Null
$d.Content.Find.Execute("469÷4=117, 1", $true, $false, $false, $false, $false, $true, 1, $false, "190÷7=27, 1", 2) | Out-Null
$d.Content.Find.Execute("361÷4=90, 1", $true, $false, $false, $false, $false, $true, 1, $false, "505÷4=126, 1", 2) | Out-Null
$d.Content.Find.Execute("255÷3=85, 0", $true, $false, $false, $false, $false, $true, 1, $false, "392÷6=65, 2", 2) | Out-Null
$d.Content.Find.Execute("471÷5=94, 1", $true, $false, $false, $false, $false, $true, 1, $false, "884÷5=176, 4", 2) | Out-Null
$d.Content.Find.Execute("418÷8=52, 2", $true, $false, $false, $false, $false, $true, 1, $false, "583÷5=116, 3", 2) | Out-Null
$d.Content.Find.Execute("939÷2=469, 1", $true, $false, $false, $false, $false, $true, 1, $false, "290÷4=72, 2", 2) | Out-Null
$d.Content.Find.Execute("617÷7=88, 1", $true, $false, $false, $false, $false, $true, 1, $false, "762÷8=95, 2", 2) | Out-Null
$d.Content.Find.Execute("757÷3=252, 1", $true, $false, $false, $false, $false, $true, 1, $false, "976÷5=195, 1", 2) | Out-Null
$d.Content.Find.Execute("917÷3=305, 2", $true, $false, $false, $false, $false, $true, 1, $false, "401÷3=133, 2", 2) | Out-Null
$d.Content.Find.Execute("895÷4=223, 3", $true, $false, $false, $false, $false, $true, 1, $false, "572÷3=190, 2", 2) | Out-Null
$d.Content.Find.Execute("637÷2=318, 1", $true, $false, $false, $false, $false, $true, 1, $false, "389÷3=129, 2", 2) | Out-Null
$d.Content.Find.Execute("165÷8=20, 5", $true, $false, $false, $false, $false, $true, 1, $false, "379÷8=47, 3", 2) | Out-Null
$d.Content.Find.Execute("941÷3=313, 2", $true, $false, $false, $false, $false, $true, 1, $false, "792÷2=396, 0", 2) | Out-Null
$d.Content.Find.Execute("574÷9=63, 7", $true, $false, $false, $false, $false, $true, 1, $false, "562÷4=140, 2", 2) | Out-Null
$d.Content.Find.Execute("702÷2=351, 0", $true, $false, $false, $false, $false, $true, 1, $false, "347÷9=38, 5", 2) | Out-Null
$d.Content.Find.Execute("473÷5=94, 3", $true, $false, $false, $false, $false, $true, 1, $false, "420÷2=210, 0", 2) | Out-Null
$d.Content.Find.Execute("420÷5=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "114÷3=38, 0", 2) | Out-Null
$d.Content.Find.Execute("890÷9=98, 8", $true, $false, $false, $false, $false, $true, 1, $false, "192÷7=27, 3", 2) | Out-Null
$d.Content.Find.Execute("806÷6=134, 2", $true, $false, $false, $false, $false, $true, 1, $false, "981÷6=163, 3", 2) | Out-Null
$d.Content.Find.Execute("993÷8=124, 1", $true, $false, $false, $false, $false, $true, 1, $false, "891÷4=222, 3", 2) | Out-Null
$d.Content.Find.Execute("408÷6=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "820÷9=91, 1", 2) | Out-Null
$d.Content.Find.Execute("947÷6=157, 5", $true, $false, $false, $false, $false, $true, 1, $false, "954÷6=159, 0", 2) | Out-Null
$d.Content.Find.Execute("582÷5=116, 2", $true, $false, $false, $false, $false, $true, 1, $false, "239÷5=47, 4", 2) | Out-Null
$d.Content.Find.Execute("767÷3=255, 2", $true, $false, $false, $false, $false, $true, 1, $false, "874÷4=218, 2", 2) | Out-Null
